# "Excel chart in different languages.xlsx" — fix the "Tabe and BackEnd"
# sheet-name typo, re-enter the (already-volatile) month formulas in
# A7:A18 as one fill so they collapse into a shared-formula group, and
# leave the user's selection/active tab on the renamed sheet at B39
# (matching the author's last-saved UI state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabe and BackEnd")

# Fix the typo in the sheet name. Excel itself keeps every formula/chart
# reference to this sheet in sync with the rename.
$ws.Name = "Table and BackEnd"

# Re-enter the same formula across A7:A18 in one shot. Typing/filling an
# identical formula across a contiguous range is what makes Excel store it
# as a single shared-formula group instead of one independent formula per
# cell.
$ws.Range("A7:A18").Formula = "=INDIRECT(Lang&ROW())"

# Make this the active sheet/selection (it became the front-most tab the
# user left selected, cursor parked on B39).
$ws.Activate()
$ws.Range("B39").Select()
